# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D3").Value = -0.288
$ws.Range("G2:G3").Value = -0.1466216216216216
$ws.Range("H2:H3").Value = -0.1466216216216216
$ws.Range("I2:I3").Value = -0.1959459459459459
$ws.Range("J2:J3").Value = -0.1959459459459459
$ws.Range("K2:K3").Value = -5.1
$ws.Range("L2:L3").Value = -0.3445945945945946
$ws.Range("U2:U3").Value = 3.26
$ws.Range("V2:V3").Value = 0.1671794871794872
$ws.Range("W2:W3").Value = -0.6166868198307134
$ws.Range("X2:X3").Value = 0.1230242329056347
$ws.Range("Y2:Y3").Value = -0.7397110527363482
$ws.Range("Z2:Z3").Value = 0.4603421461897357
$ws.Range("AA2:AA3").Value = -0.09020217729393468
$ws.Range("AB2:AB3").Value = 0.06449127636040566
$ws.Range("AC2:AC3").Value = -0.1546934536543403
$ws.Range("AD2:AD3").Value = 29.5
$ws.Range("AE2:AE3").Value = 0
$ws.Range("AF2:AF3").Value = 29.5
$ws.Range("AG2:AG3").Value = 26.24
$ws.Range("AH2:AH3").Value = 0.6020408163265306
$ws.Range("AI2:AI3").Value = 0.9382951653944019
$ws.Range("AJ2:AJ3").Value = 0.5736773065150853
$ws.Range("AK2:AK3").Value = 0.9311568488289567
$ws.Range("AL2:AL3").Value = 0.958
$ws.Range("AM2:AM3").Value = 0.82
$ws.Range("AN2:AN3").Value = -11.75298804780877
$ws.Range("AO2:AO3").Value = -3.02713987473904
$ws.Range("AP2:AP3").Value = -10.45418326693227
$ws.Range("AQ2:AQ3").Value = -3.536585365853659
